$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (GitHub Actions scheduled update)

$ws.Range('D2').Value = '41.815.84'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '2.264.53'
$ws.Range('E3').Value = '  -0.25%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'303.51"
$ws.Range('E5').Value = '  +0.56%  '
$ws.Range('D6').Value = "'92.63"
$ws.Range('E6').Value = '  -0.20%  '
$ws.Range('E7').Value = '  +0.83%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -0.62%  '
$ws.Range('D10').Value = "'32.41"
$ws.Range('E10').Value = '  -0.74%  '
$ws.Range('E11').Value = '  -0.33%  '
$ws.Range('E12').Value = '  -1.91%  '
$ws.Range('D13').Value = "'6.65"
$ws.Range('E13').Value = '  -0.73%  '
$ws.Range('D14').Value = '2.615.88'
$ws.Range('E14').Value = '  -0.21%  '
$ws.Range('D15').Value = "'14.29"
$ws.Range('E15').Value = '  +0.73%  '
$ws.Range('D16').Value = '2.264.18'
$ws.Range('E16').Value = '  -0.74%  '
$ws.Range('D17').Value = "'0.783"
$ws.Range('E17').Value = '  +3.71%  '
$ws.Range('D18').Value = '41.744.48'
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('D19').Value = "'12.73"
$ws.Range('E19').Value = '  +2.83%  '
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').Value = "'5.92"
$ws.Range('E21').Value = '  -0.23%  '
$ws.Range('D22').Value = "'67.63"
$ws.Range('E22').Value = '  +0.59%  '
$ws.Range('D23').Value = "'244.48"
$ws.Range('E23').Value = '  +1.40%  '
$ws.Range('E24').Value = '  +0.52%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('E26').Value = '  +1.52%  '
$ws.Range('D27').Value = "'24.02"
$ws.Range('E27').Value = '  +0.87%  '
$ws.Range('E28').Value = '  -1.25%  '
$ws.Range('E29').Value = '  -5.73%  '
$ws.Range('D30').Value = "'34.97"
$ws.Range('E30').Value = '  +2.67%  '
$ws.Range('D31').Value = "'158.88"
$ws.Range('E31').Value = '  +0.50%  '
$ws.Range('E32').Value = '  +1.84%  '
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('E34').Value = '  +0.57%  '
$ws.Range('D35').Value = "'3.02"
$ws.Range('E35').Value = '  -1.52%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').Value = "'0.106"
$ws.Range('E36').Value = '  +1.69%  '
$ws.Range('B37').Value = 'Celestia'
$ws.Range('C37').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D37').Value = "'16.85"
$ws.Range('E37').Value = '  +2.12%  '
$ws.Range('E38').Value = '  -1.04%  '
$ws.Range('E39').Value = '  +0.84%  '
$ws.Range('D40').Value = "'1.79"
$ws.Range('E40').Value = '  +0.50%  '
$ws.Range('E41').Value = '  -1.65%  '
$ws.Range('D42').Value = "'20.09"
$ws.Range('E42').Value = '  -0.25%  '
$ws.Range('D43').Value = '2.006.53'
$ws.Range('E43').Value = '  -2.00%  '
$ws.Range('D44').Value = "'2.25"
$ws.Range('E44').Value = '  +12.51%  '
$ws.Range('E45').Value = '  +1.43%  '
$ws.Range('D46').Value = "'10.37"
$ws.Range('E46').Value = '  +2.67%  '
$ws.Range('E47').Value = '  -1.56%  '
$ws.Range('D48').Value = "'73.16"
$ws.Range('E48').Value = '  +3.67%  '
$ws.Range('D49').Value = "'52.88"
$ws.Range('E49').Value = '  +2.34%  '
$ws.Range('E50').Value = '  +0.78%  '
$ws.Range('D51').Value = "'1.51"
$ws.Range('E51').Value = '  -0.10%  '

Write-Output "Applied cryptos update"
